$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C new cell values (translated text) ---
$ws.Range("C5").Value = "\n<\n[3]>Smooch♥ Sluuuurp♥ Hoora♥`nYou wanted some of this, that's why you were peeking, right?`nLick, lick♥"
$ws.Range("C6").Value = "\n<\n[3]>If you want to be raped so badly,`njust be honest and tell me.`nSluuuurp♥ Lick...♥"
$ws.Range("C7").Value = "\n<\n[3]>Hora♥ Lick♥`nWhen I do this, you won't be able to run away, right?`nYou don't have to run away. Sluurp♥"
$ws.Range("C8").Value = "\n<\n[3]>Caught you♥`nSluuuuurp♥ Lick♥ Smoooch♥`nAhaaa♥ I taste some precum♥"
$ws.Range("C9").Value = "\n<\n[3]>Hmmm?`nAre you gonna cum if I keep licking?`nShould I stop? Lick, lick, lick♥"
$ws.Range("C10").Value = "\n<\n[3]>Mmm♥`nYou came...♥`nSluuuuuuurp♥"
$ws.Range("C12").Value = "Nothing else important."
$ws.Range("C16").Value = "\n<\n[1]>Hot!!!"
$ws.Range("C22").Value = "This is bad…! Gotta get out…!"
$ws.Range("C23").Value = "Start running"
$ws.Range("C24").Value = "Hint"
$ws.Range("C25").Value = "Head for the exit as Lily and Shina chase you.`nLime will capture you if you step in a puddle.`nIf your energy is low, use an item or the Deep Breathing skill."
$ws.Range("C26").Value = "Check strategy"
$ws.Range("C27").Value = "That's enough"
$ws.Range("C28").Value = "The puddles disappear after being stepped on.`nIf you think you can't avoid them, just push on through."
$ws.Range("C32").Value = "You managed to escape…"
$ws.Range("C34").Value = "\n<\n[3]>So you came here to get stepped on-?`nStomp stomp stomp-♥`nNyahahahaha♥"
$ws.Range("C35").Value = "\n<\n[3]>The punishment for peeping is a footjob, nyan♥`nGrind, grind♥`nWhat's that, nya? You want me to step on you some more-?"
$ws.Range("C36").Value = "\n<\n[3]>Nyaha♥`nYou submit the moment you get stepped on♥`nThat's your natural masochistic instinct, nyan♥"
$ws.Range("C37").Value = "\n<\n[3]>Pretty ballsy of you to interrupt my bath, nyan.`nI'll squeeze out some sperm to use for bath salts♥`nNyahahahaha!"
$ws.Range("C38").Value = "\n<\n[3]>Haaah～?`nYou came already, nyaa?`nThen it's time for the killing blow, nya♥"
$ws.Range("C39").Value = "\n<\n[3]>Stomp stomp stomp stomp-♥"
$ws.Range("C40").Value = "\n<\n[3]>Ahhh. You came so much, nya.`nBut I'm not done stepping on your, nyaa.`nAre you finished being stepped on-?"
$ws.Range("C45").Value = "\>\C[26]Avoided the trap!"
$ws.Range("C46").Value = "\>\C[14]Stepped on a spring trap!\C[0]"

# --- Wrap text for long dialogue / hint cells ---
$ws.Range("B5:C5").WrapText = $true
$ws.Range("B6:C6").WrapText = $true
$ws.Range("B7:C7").WrapText = $true
$ws.Range("B8:C8").WrapText = $true
$ws.Range("B9:C9").WrapText = $true
$ws.Range("B10:C10").WrapText = $true
$ws.Range("A17").WrapText = $true
$ws.Range("B25:C25").WrapText = $true
$ws.Range("B28:C28").WrapText = $true
$ws.Range("B34:C34").WrapText = $true
$ws.Range("B35:C35").WrapText = $true
$ws.Range("B36:C36").WrapText = $true
$ws.Range("B37:C37").WrapText = $true
$ws.Range("B38:C38").WrapText = $true
$ws.Range("B40:C40").WrapText = $true

# --- Row heights for wrapped rows ---
$ws.Rows.Item(5).RowHeight = 45
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 45
$ws.Rows.Item(10).RowHeight = 45
$ws.Rows.Item(17).RowHeight = 75
$ws.Rows.Item(25).RowHeight = 45
$ws.Rows.Item(28).RowHeight = 30
$ws.Rows.Item(34).RowHeight = 45
$ws.Rows.Item(35).RowHeight = 45
$ws.Rows.Item(36).RowHeight = 45
$ws.Rows.Item(37).RowHeight = 45
$ws.Rows.Item(38).RowHeight = 45
$ws.Rows.Item(40).RowHeight = 45

# --- Column widths (A, B, C) ---
$ws.Columns.Item(1).ColumnWidth = 56.0
$ws.Columns.Item(2).ColumnWidth = 71.16666666666667
$ws.Columns.Item(3).ColumnWidth = 67.83333333333333

Write-Host "Edit complete"
